# Added new test case ID 13
$wb = $excel.ActiveWorkbook

$runManager = $wb.Worksheets.Item("RUNMANAGER")
$adminFunctionality = $wb.Worksheets.Item("ADMINFUNCTIONALITY")

# --- RUNMANAGER: add the new test case row (row 12) ---
$runManager.Cells.Item(12, 1).Value = "verifyThaTheNumberOfUserRecordsIncreasesAfterAddingANewUser"
$runManager.Cells.Item(12, 2).Value = "To check this test is executed"
$runManager.Cells.Item(12, 3).Value = "yes"
$runManager.Cells.Item(12, 4).Value = "'11"
$runManager.Cells.Item(12, 5).Value = "'1"

# --- ADMINFUNCTIONALITY: update employeename test data & append two rows ---
$adminFunctionality.Cells.Item(2, 5).Value = "Orange Middle NameMiddle Name Test"
$adminFunctionality.Cells.Item(3, 5).Value = "Orange Middle NameMiddle Name Test"

$adminFunctionality.Cells.Item(4, 1).Value = "verifyThatTheAdminCanAddNewUser"
$adminFunctionality.Cells.Item(4, 2).Value = "yes"
$adminFunctionality.Cells.Item(4, 3).Value = "Admin"
$adminFunctionality.Cells.Item(4, 4).Value = "admin123"
$adminFunctionality.Cells.Item(4, 5).Value = "Orange Middle NameMiddle Name Test"
$adminFunctionality.Cells.Item(4, 6).Value = "chrome"
$adminFunctionality.Cells.Item(4, 7).Value = "sham12345"
$adminFunctionality.Cells.Item(4, 8).Value = "sham12345"

$adminFunctionality.Cells.Item(5, 1).Value = "verifyThatTheAdminCanAddNewUser"
$adminFunctionality.Cells.Item(5, 2).Value = "yes"
$adminFunctionality.Cells.Item(5, 3).Value = "Admin"
$adminFunctionality.Cells.Item(5, 4).Value = "admin123"
$adminFunctionality.Cells.Item(5, 5).Value = "Orange Middle NameMiddle Name Test"
$adminFunctionality.Cells.Item(5, 6).Value = "firefox"
$adminFunctionality.Cells.Item(5, 7).Value = "sham12345"
$adminFunctionality.Cells.Item(5, 8).Value = "sham12345"

# --- Column width adjustments on ADMINFUNCTIONALITY to fit the new longer text ---
# (widths are re-computed the way Excel's own "best fit" autosize would do it for
# the new, longer cell contents; the values below are chosen so that after the
# host's column-width rounding they land as close as possible to Excel's real
# best-fit pixel widths of 54.7890625 / 32.1015625 / 14.41796875 characters)
$adminFunctionality.Columns.Item(1).ColumnWidth = 54
$adminFunctionality.Columns.Item(5).ColumnWidth = 31.3333333
$adminFunctionality.Columns.Item(8).ColumnWidth = 13.6666667

# --- Selection / active sheet adjustments ---
$adminFunctionality.Range("E2").Select() | Out-Null
$runManager.Select() | Out-Null
$runManager.Range("D12").Select() | Out-Null
